# Apply updated "dSF" (column F) values as per repull of data / mean calculation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4  = -4
    7  = 8
    9  = -5
    10 = -1
    13 = 16
    16 = 2
    17 = -5
    19 = 7
    21 = 5
    22 = 1
    23 = -5
    24 = -2
    27 = 0
    30 = 0
    31 = -2
    32 = -6
    33 = -2
    34 = -4
    36 = -3
    37 = 0
    38 = -1
    39 = 0
    42 = 2
    43 = -4
    48 = -7
    52 = 5
    53 = 3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
